# se actualizo la Clase de Equivalencia para Categoria
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clase equivalencia categoria")

# Fill in the new "Descripcion" equivalence-class block (rows 8-10),
# mirroring the existing "Identificador" block in rows 5-7.
$ws.Range("C8").Value = "Descripcion"
$ws.Range("D8").Value = "Logico"
$ws.Range("E8").Value = "descripcion= caracteres alfanumericos"
$ws.Range("F8").Value = "CEV<2>"
$ws.Range("G8").Value = "descripcion != caracteres alfanumericos"
$ws.Range("H8").Value = "CENV<04>"

$ws.Range("D9").Value = "Valor"
$ws.Range("E9").Value = "0<=ID<=50"
$ws.Range("F9").Value = "CEV<2>"
$ws.Range("G9").Value = "descripcion<1"
$ws.Range("H9").Value = "CENV<05>"

$ws.Range("G10").Value = "descripcion>50"
$ws.Range("H10").Value = "CENV<06>"

# Move the active selection the way it ended up after the edit.
$ws.Range("H16").Select() | Out-Null
